$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.720.10'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').Value = '2.624.67'
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'519.75"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.41%  '
$ws.Range('D6').Value = "'144.23"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = "'0.568"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.83%  '
$ws.Range('D9').Value = '2.631.24'
$ws.Range('E9').Value = '  +1.56%  '
$ws.Range('E10').Value = '  -0.05%  '
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('E12').Value = '  -0.75%  '
$ws.Range('E13').Value = '  -0.96%  '
$ws.Range('D14').Value = '3.084.60'
$ws.Range('E14').Value = '  +1.97%  '
$ws.Range('D15').Value = '58.742.82'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').Value = "'20.72"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.67%  '
$ws.Range('E17').Value = '  -0.87%  '
$ws.Range('D18').Value = '2.630.30'
$ws.Range('E18').Value = '  +1.97%  '
$ws.Range('D19').Value = "'344.86"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.69%  '
$ws.Range('D20').Value = "'4.44"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.21%  '
$ws.Range('D21').Value = "'10.16"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').Value = "'6.12"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.04%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  +1.41%  '
$ws.Range('E25').Value = '  -1.27%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').Value = "'0.997"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = "'0.163"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.16%  '
$ws.Range('D28').Value = '0.0₃0795'
$ws.Range('E28').Value = '  -2.24%  '
$ws.Range('E29').Value = '  +0.80%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('D31').Value = "'6.21"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.65%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = "'18.81"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.11%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = "'1.57"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.98%  '
$ws.Range('D34').Value = "'150.00"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('D35').Value = "'0.975"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.76%  '
$ws.Range('E36').Value = '  -0.30%  '
$ws.Range('D37').Value = "'1.13"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('D38').Value = "'36.62"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.63%  '
$ws.Range('D39').Value = "'0.834"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.00%  '
$ws.Range('D40').Value = "'3.63"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.29%  '
$ws.Range('E41').Value = '  +1.36%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = "'0.997"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').Value = "'276.85"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.35%  '
$ws.Range('D44').Value = "'0.0980"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.98%  '
$ws.Range('D45').Value = "'0.598"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.88%  '
$ws.Range('D46').Value = "'19.46"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.11%  '
$ws.Range('D47').Value = "'10.31"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('D48').Value = "'0.0519"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.11%  '
$ws.Range('D49').Value = '1.987.79'
$ws.Range('E49').Value = '  +3.21%  '
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('D51').Value = "'4.62"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.26%  '

Write-Output "Applied 97 cell updates"
